{"js": "// Replace the long, itemised \"Micro results\" history in the second\n// (answer) cell of the \"Micro results\" row with a short blue note that\n// archives the old results behind a \"Previous result (1 year)\" divider,\n// matching the author's \"update micro results summrization\" edit.\n\n// Locate the results table and the \"Micro results\" row by scanning the\n// first cell of each row for the label text.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nlet targetCell = null;\nfor (const row of rows.items) {\n  const cells = row.cells.items;\n  if (cells.length >= 2 && cells[0].body.text.trim() === \"Micro results\") {\n    targetCell = cells[1];\n    break;\n  }\n}\n\nif (!targetCell) {\n  throw new Error('Could not find the \"Micro results\" row.');\n}\n\n// Flat-OPC wrapped WordprocessingML fragment for the three paragraphs\n// that replace the long culture-result history: a blank (formatted but\n// empty) run, the \"Previous result\" divider line, and a trailing blank\n// (formatted but empty) run.\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:color w:val=\"0000FF\"/><w:sz w:val=\"20\"/></w:rPr></w:r></w:p>' +\n  '<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:color w:val=\"0000FF\"/><w:sz w:val=\"20\"/></w:rPr><w:t>--------Previous result (1 year)--------</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:color w:val=\"0000FF\"/><w:sz w:val=\"20\"/></w:rPr></w:r></w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// Clearing the cell body leaves a single empty paragraph (matching the\n// diff's leading `<w:p><w:r/></w:p>`); insert the remaining three\n// paragraphs after it.\ntargetCell.body.clear();\ntargetCell.body.insertOoxml(replacementOoxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Replace the long, itemised \"Micro results\" history in the second\n# (answer) cell of the \"Micro results\" row with a short blue note that\n# archives the old results behind a \"Previous result (1 year)\" divider,\n# matching the author's \"update micro results summrization\" edit.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Locate the \"Micro results\" row by scanning the first cell of each row\n# (cell text carries a trailing cell-mark, so trim it before comparing).\n$targetRow = 0\nfor ($i = 1; $i -le $table.Rows.Count; $i++) {\n    $label = $table.Cell($i, 1).Range.Text.Trim([char]13, [char]7)\n    if ($label -eq \"Micro results\") {\n        $targetRow = $i\n        break\n    }\n}\n\nif ($targetRow -eq 0) {\n    throw \"Could not find the 'Micro results' row.\"\n}\n\n$cell = $table.Cell($targetRow, 2)\n\n# Flat-OPC wrapped WordprocessingML fragment for the four paragraphs that\n# replace the long culture-result history: an empty paragraph, a blank\n# (formatted but empty) run, the \"Previous result\" divider line, and a\n# trailing blank (formatted but empty) run. InsertXML replaces the whole\n# contents of the target range.\n$xml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n</Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:r/></w:p>\n<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:color w:val=\"0000FF\"/><w:sz w:val=\"20\"/></w:rPr></w:r></w:p>\n<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:color w:val=\"0000FF\"/><w:sz w:val=\"20\"/></w:rPr><w:t>--------Previous result (1 year)--------</w:t></w:r></w:p>\n<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:color w:val=\"0000FF\"/><w:sz w:val=\"20\"/></w:rPr></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$cell.Range.InsertXML($xml)\n"}
